# New weekly price record for "Espinaca" (Vega Modelo de Temuco) needs to be
# inserted as the new first record of that block (row 92), pushing the
# existing historical rows (old 92-103) down by one row (new 93-104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 92 (shifts rows 92:103 down to 93:104)
$ws.Rows.Item(92).Insert()

# Populate the new row 92 by duplicating the row that just got pushed down
# to row 93 (this preserves formatting/styles and all the columns that stay
# constant for this data block), then overwrite the two columns that hold
# the new week's actual values (Fecha / Volumen).
$ws.Range("A93:R93").Copy($ws.Range("A92:R92"))

$ws.Cells.Item(92, 4).Value2 = 44522   # D92 - Fecha
$ws.Cells.Item(92, 10).Value2 = 50     # J92 - Volumen
